$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the order status for the first two orders to "Cancelled"
$ws.Range("E2").Value = "Cancelled"
$ws.Range("E3").Value = "Cancelled"
